# tradexcb_strategy.xlsx — update the sample trade row on "Sheet1":
#   * flip the sample order from a Buy/LIMIT entry to a Sell/MARKET entry
#   * roll the expiry/instrument to the next weekly contract
#   * bump lot size (stoploss), tsl and target figures, and the timeframe
# (matches: "Updated time profiling code + Fixed slicing issue, where
#  keeps on placing orders when n_slices > no. of lots")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Plain text columns - safe to assign directly.
$ws.Range("A2").Value = "Sell"
$ws.Range("D2").Value = "MARKET"
$ws.Range("J2").Value = "NIFTY2250517000CE"

# Columns that look numeric/date-like must stay plain text (as the sheet
# stores every value as text). A leading apostrophe forces text entry the
# same way typing it into Excel would, and re-applying the "Normal" cell
# style afterwards clears the resulting quote-prefix formatting so the
# cell's style stays the same as before the edit.
$ws.Range("I2").Value = "'2022-05-05"
$ws.Range("L2").Value = "'20"
$ws.Range("N2").Value = "'5"
$ws.Range("P2").Value = "'40"
$ws.Range("Q2").Value = "'4"

$ws.Range("I2").Style = "Normal"
$ws.Range("L2").Style = "Normal"
$ws.Range("N2").Style = "Normal"
$ws.Range("P2").Style = "Normal"
$ws.Range("Q2").Style = "Normal"

# The author left the workbook with "Sheet1" active and the selection
# parked on M15.
$ws.Activate()
$ws.Range("M15").Select()
